$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new criteria row "Kosten" above the old row 9 ("Total Punkte") ---
# This pushes "Total Punkte" (old row 9) down to row 10, "Rang" (old row 10) to
# row 11, the blank spacer (old row 11) to row 12, and the two remark rows
# (old rows 12/13) down to rows 13/14.
$ws.Rows(9).Insert()

# Copy the formatting (styles, borders, shared-formula cell styles) of the row
# above ("Darstellungsqualitaet...", now row 8) into the freshly inserted row 9.
$ws.Range("A8:H8").Copy($ws.Range("A9:H9"))

# Fill in the new "Kosten" criterion row.
$ws.Range("A9").Value = "Kosten"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Formula = "=B9*C9"
$ws.Range("E9").Value = 5
$ws.Range("F9").Formula = "=B9*E9"
$ws.Range("G9").Value = 3
$ws.Range("H9").Formula = "=B9*G9"

# --- Update the wording of the first criterion row (row 5) ---
$ws.Range("A5").Value = "Einbringung in Raumhöhe (Raumgefühl)"
$ws.Range("E5").Value = 1

# --- Update scores on "Bildschirmfläche überblickbar" (row 6) ---
$ws.Range("C6").Value = 5
$ws.Range("G6").Value = 3

# --- Update scores on "Darstellungsqualität/-grösse zum Lesen" (row 8) ---
$ws.Range("G8").Value = 1

# --- Fix the "Total Punkte" sums (now row 10) to include the new row 9 ---
$ws.Range("D10").Formula = "=SUM(D5:D9)"
$ws.Range("F10").Formula = "=SUM(F5:F9)"
$ws.Range("H10").Formula = "=SUM(H5:H9)"

# --- Fix the "Rang" formulas (now row 11) to reference the shifted totals row ---
$ws.Range("D11").Formula = "=IF(D10>=F10,IF(D10>=H10,1,2),IF(D10>=H10,2,3))"
$ws.Range("F11").Formula = "=IF(F10>=H10,IF(F10>=D10,1,2),IF(F10>=D10,2,3))"
$ws.Range("H11").Formula = "=IF(H10>=D10,IF(H10>=F10,1,2),IF(H10>=F10,2,3))"

# --- Update the remark text (now row 13) with the new wording ---
$ws.Range("A13").Value = "Bemerkung: Die Gewichtungs- / Bewertungsskala geht von wenig (1), bedingt (3) bis zu sehr wichtig (5)."

# --- Update the selected cell to match the saved state ---
$ws.Range("H15").Select() | Out-Null
